$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the 39 new Q&A rows (rows 38-76) coming from the Gemini chatbot update
$ws.Cells.Item(38, 1).Value = 'What are the timings of sports complex?'
$ws.Cells.Item(38, 2).Value = '6 am to 8 am and 4:45 pm to 7:00 pm'
$ws.Cells.Item(39, 1).Value = 'what are the swiming timings for men?'
$ws.Cells.Item(39, 2).Value = '6am-7am and 5pm to 6pm'
$ws.Cells.Item(40, 1).Value = 'what are the swiming timings for women?'
$ws.Cells.Item(40, 2).Value = '7am-8am and 6pm to 7pm'
$ws.Cells.Item(41, 1).Value = 'what are hostels visiting timings?'
$ws.Cells.Item(41, 2).Value = '10 am to 5pm'
$ws.Cells.Item(42, 1).Value = 'hostel timings ?'
$ws.Cells.Item(42, 2).Value = '6am to 7pm'
$ws.Cells.Item(43, 1).Value = 'what is library timings?'
$ws.Cells.Item(43, 2).Value = 'All Working Days : 9:00Am to 11:00 PM, Second Saturday & Sundays : 09:30 Am to 4:30 PM and closed on Public Holidays'
$ws.Cells.Item(44, 1).Value = 'From where should we receive memos?'
$ws.Cells.Item(44, 2).Value = 'In SIT building'
$ws.Cells.Item(45, 1).Value = 'About JNTUH college?'
$ws.Cells.Item(45, 2).Value = 'Jawaharlal Nehru Technological University College of Engineering Hyderabad (Autonomous), formerly known as Nagarjuna Sagar Engineering College, was established in 1965 by the Government of Andhra Pradesh and administrated under the control of the Department of Technical Education and affiliated to Osmania University, Hyderabad. With the formation of Jawaharlal Nehru Technological University(JNTU) on 2nd October 1972, the college was made as a constituent college of the University and eventually renamed as JNTU College of Engineering, Hyderabad.'
$ws.Cells.Item(46, 1).Value = 'What are the documents required to submit while admissions?'
$ws.Cells.Item(46, 2).Value = '10th class memo, inter memo, caste certificate, aadhar Xerox,eamcet rank card, allotment letter, joining letter'
$ws.Cells.Item(47, 1).Value = 'How can I contact campus security?'
$ws.Cells.Item(47, 2).Value = 'Dial the security hotline number available on your ID card.'
$ws.Cells.Item(48, 1).Value = 'What is the fee for late book returns?'
$ws.Cells.Item(48, 2).Value = 'The late fee is $0.50 per day per book.'
$ws.Cells.Item(49, 1).Value = 'What is the cost of a meal in the cafeteria?'
$ws.Cells.Item(49, 2).Value = 'The average meal costs around $5 to $7.'
$ws.Cells.Item(50, 1).Value = 'Where can I park my bike?'
$ws.Cells.Item(50, 2).Value = 'There is a bike parking area near the main entrance.'
$ws.Cells.Item(51, 1).Value = 'How do I join the student council?'
$ws.Cells.Item(51, 2).Value = 'Submit an application during the student council elections.'
$ws.Cells.Item(52, 1).Value = 'Where can I get a transcript of my grades?'
$ws.Cells.Item(52, 2).Value = 'Request transcripts at the registrarâ€™s office.'
$ws.Cells.Item(53, 1).Value = 'Are there lockers available for students?'
$ws.Cells.Item(53, 2).Value = 'Yes, lockers can be rented from the administration office.'
$ws.Cells.Item(54, 1).Value = 'Are there any vegan options in the cafeteria?'
$ws.Cells.Item(54, 2).Value = 'Yes, vegan meals are available on request.'
$ws.Cells.Item(55, 1).Value = 'What is the procedure to apply for a library card?'
$ws.Cells.Item(55, 2).Value = 'Fill out the application form at the library reception.'
$ws.Cells.Item(56, 1).Value = 'Are pets allowed on campus?'
$ws.Cells.Item(56, 2).Value = 'No, pets are not permitted on campus grounds.'
$ws.Cells.Item(57, 1).Value = 'What is the deadline to apply for scholarships?'
$ws.Cells.Item(57, 2).Value = 'Scholarship deadlines vary; check the website for details.'
$ws.Cells.Item(58, 1).Value = 'How do I apply for a research grant?'
$ws.Cells.Item(58, 2).Value = 'Submit a proposal through the university''s research office.'
$ws.Cells.Item(59, 1).Value = 'Can visitors use the library facilities?'
$ws.Cells.Item(59, 2).Value = 'Yes, visitors can use library facilities with a guest pass.'
$ws.Cells.Item(60, 1).Value = 'Where can I find information on upcoming events?'
$ws.Cells.Item(60, 2).Value = 'Events are posted on the bulletin board and the university website.'
$ws.Cells.Item(61, 1).Value = 'How can I get a duplicate ID card?'
$ws.Cells.Item(61, 2).Value = 'Apply for a duplicate at the administration office with a fee.'
$ws.Cells.Item(62, 1).Value = 'What are the gym opening hours?'
$ws.Cells.Item(62, 2).Value = 'The gym is open from 6 am to 10 pm daily.'
$ws.Cells.Item(63, 1).Value = 'How do I access Wi-Fi on campus?'
$ws.Cells.Item(63, 2).Value = 'Login with your student credentials; guest access is also available.'
$ws.Cells.Item(64, 1).Value = 'How do I change my course schedule?'
$ws.Cells.Item(64, 2).Value = 'Contact the registrarâ€™s office for schedule changes.'
$ws.Cells.Item(65, 1).Value = 'What is the dress code for lab classes?'
$ws.Cells.Item(65, 2).Value = 'Lab coats and closed-toe shoes are required in lab classes.'
$ws.Cells.Item(66, 1).Value = 'Is there a lost and found service?'
$ws.Cells.Item(66, 2).Value = 'Yes, lost items are handed to the security desk at the main gate.'
$ws.Cells.Item(67, 1).Value = 'Where can I find first aid on campus?'
$ws.Cells.Item(67, 2).Value = 'Visit the campus health center or contact security for assistance.'
$ws.Cells.Item(68, 1).Value = 'Are there part-time job opportunities on campus?'
$ws.Cells.Item(68, 2).Value = 'Yes, check the job board or speak to the career center.'
$ws.Cells.Item(69, 1).Value = 'Are laptops allowed in the library?'
$ws.Cells.Item(69, 2).Value = 'Yes, laptops are allowed in designated areas of the library.'
$ws.Cells.Item(70, 1).Value = 'Is there a shuttle service on campus?'
$ws.Cells.Item(70, 2).Value = 'Yes, the campus shuttle operates from 8 am to 8 pm.'
$ws.Cells.Item(71, 1).Value = 'Is smoking allowed on campus?'
$ws.Cells.Item(71, 2).Value = 'No, smoking is prohibited in all campus areas.'
$ws.Cells.Item(72, 1).Value = 'How do I access the online learning portal?'
$ws.Cells.Item(72, 2).Value = 'Login with your student ID and password on the website.'
$ws.Cells.Item(73, 1).Value = 'What is the process for joining a club?'
$ws.Cells.Item(73, 2).Value = 'Attend the club orientation or register during club day.'
$ws.Cells.Item(74, 1).Value = 'Can I book a study room in the library?'
$ws.Cells.Item(74, 2).Value = 'Yes, study rooms can be booked at the library reception.'
$ws.Cells.Item(75, 1).Value = 'How can I reset my email password?'
$ws.Cells.Item(75, 2).Value = 'Visit the IT support office or reset through the university portal.'
$ws.Cells.Item(76, 1).Value = 'Where can I find help with my assignments?'
$ws.Cells.Item(76, 2).Value = 'Visit the academic support center or ask your professor.'

# Restore the view state recorded in the workbook (selection + zoom)
[void]$ws.Cells.Item(61, 7).Select()
$excel.ActiveWindow.Zoom = 51
